# [RF GEN]CPU SCH update
# Populate Sheet1 with the CPU SCH RF generator review data table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Text labels/units first, in authoring order (drives shared-string order) ---
$ws.Range("C5").Value = "VIN"
$ws.Range("G5").Value = "V"
$ws.Range("C6").Value = "R1"
$ws.Range("C7").Value = "R2"
$ws.Range("G6").Value = "KΩ"
$ws.Range("G7").Value = "KΩ"
$ws.Range("C8").Value = "I_IN"
$ws.Range("G8").Value = "mA"
$ws.Range("C9").Value = "P_IN"
$ws.Range("G9").Value = "W"
$ws.Range("C10").Value = "V_R2"
$ws.Range("G10").Value = "V"

# --- Row 5: VIN values ---
$ws.Range("D5").Value = 48
$ws.Range("E5").Formula = "=48*0.8"
$ws.Range("F5").Formula = "=48*1.2"

# --- Row 6: R1 values ---
$ws.Range("D6").Value = 20
$ws.Range("E6").Value = 20
$ws.Range("F6").Value = 20

# --- Row 7: R2 values ---
$ws.Range("D7").Value = 4.99
$ws.Range("E7").Value = 4.99
$ws.Range("F7").Value = 4.99

# --- Row 8: I_IN values ---
$ws.Range("D8").Formula = "=D5/(D6+D7)"
$ws.Range("E8").Formula = "=E5/(E6+E7)"
$ws.Range("F8").Formula = "=F5/(F6+F7)"

# --- Row 9: P_IN values ---
$ws.Range("D9").Formula = "=D5*D8/1000"
$ws.Range("E9").Formula = "=E5*E8/1000"
$ws.Range("F9").Formula = "=F5*F8/1000"

# --- Row 10: V_R2 values ---
$ws.Range("D10").Formula = "=D5*D7/(D6+D7)"
$ws.Range("E10").Formula = "=E5*E7/(E6+E7)"
$ws.Range("F10").Formula = "=F5*F7/(F6+F7)"

# The computed ratio/power/voltage rows carry a 2-decimal-place number format.
$ws.Range("D8:F10").NumberFormat = "0.00"

# Page setup: A4, portrait.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$ws.Range("J20").Select()
